# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (strikeouts -> K count), per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 4
    4  = 3
    5  = 7
    6  = 2
    7  = 4
    8  = 2
    9  = 10
    10 = 3
    11 = 4
    12 = 1
    13 = 4
    14 = 0
    15 = 2
    16 = 6
    17 = 2
    18 = 6
    19 = 2
    20 = 4
    21 = 6
    22 = 2
    23 = 4
    24 = 6
    25 = 4
    26 = 3
    27 = 3
    28 = 2
    29 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
